$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 104 (after update) ---
$ws.Range("B104").Value = 7127370
$ws.Range("F104").Value = "Macarthur FC"
$ws.Range("G104").Value = "Wellington Phoenix"
$ws.Range("I104").Value = 2
$ws.Range("J104").Value = "A"
$ws.Range("K104").Value = 2.4
$ws.Range("M104").Value = 2.625
$ws.Range("N104").Value = 2.375
$ws.Range("O104").Value = 3.8
$ws.Range("P104").Value = 2.75
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = 1.8
$ws.Range("S104").Value = 2.05
$ws.Range("T104").Value = 3
$ws.Range("U104").Value = 1.9
$ws.Range("V104").Value = 1.95
$ws.Range("W104").Value = -1
$ws.Range("Y104").Value = 1.75
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = 1.05
$ws.Range("AB104").Value = 0
$ws.Range("AC104").Value = -0

# --- Row 105 (after update) ---
$ws.Range("B105").Value = 7127374
$ws.Range("F105").Value = "Central Coast Mariners"
$ws.Range("G105").Value = "Western Sydney Wanderers"
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = "H"
$ws.Range("K105").Value = 1.909
$ws.Range("M105").Value = 3.6
$ws.Range("N105").Value = 2.15
$ws.Range("O105").Value = 3.6
$ws.Range("P105").Value = 3.25
$ws.Range("Q105").Value = -0.25
$ws.Range("R105").Value = 1.86
$ws.Range("S105").Value = 2.04
$ws.Range("T105").Value = 2.75
$ws.Range("U105").Value = 1.975
$ws.Range("V105").Value = 1.875
$ws.Range("W105").Value = 1.15
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = 0.8600000000000001
$ws.Range("AA105").Value = -1
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0.875

# --- Row 146 ---
$ws.Range("R146").Value = 1.84
$ws.Range("S146").Value = 2.06

# --- Row 147 ---
$ws.Range("N147").Value = 1.75
$ws.Range("P147").Value = 4
$ws.Range("R147").Value = 1.97
$ws.Range("S147").Value = 1.93
$ws.Range("T147").Value = 3.75
$ws.Range("U147").Value = 1.95
$ws.Range("V147").Value = 1.9

# --- Row 148 ---
$ws.Range("N148").Value = 3.6
$ws.Range("R148").Value = 1.9
$ws.Range("S148").Value = 2
$ws.Range("U148").Value = 1.85
$ws.Range("V148").Value = 2

# --- Row 149 ---
$ws.Range("N149").Value = 1.727
$ws.Range("O149").Value = 4.333
$ws.Range("P149").Value = 4.2
$ws.Range("R149").Value = 1.9
$ws.Range("S149").Value = 2

# --- Row 150 ---
$ws.Range("O150").Value = 4.5
$ws.Range("P150").Value = 4.75
$ws.Range("R150").Value = 2.06
$ws.Range("S150").Value = 1.84

# --- Row 151 ---
$ws.Range("P151").Value = 7
$ws.Range("R151").Value = 1.88
$ws.Range("S151").Value = 2.02
$ws.Range("U151").Value = 2.025
$ws.Range("V151").Value = 1.825
